$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0009840167343248343
$ws.Range("E2").Value = 0.0009840167343248343

# Row 3
$ws.Range("D3").Value = 0.4923721784310911
$ws.Range("E3").Value = 0.4923721784310911

# Row 4
$ws.Range("D4").Value = 0.003552574232653041
$ws.Range("E4").Value = 0.003552574232653041

# Row 5
$ws.Range("D5").Value = 0.001020292602139061
$ws.Range("E5").Value = 0.001020292602139061

# Row 6
$ws.Range("D6").Value = 0.2672365802257742
$ws.Range("E6").Value = 0.2672365802257742

# Row 7
$ws.Range("D7").Value = 0.9999999999972036
$ws.Range("E7").Value = 2.796429754425844 / 1000000000000

# Row 8
$ws.Range("D8").Value = 0.9734963279178319
$ws.Range("E8").Value = 0.02650367208216808

# Row 9
$ws.Range("D9").Value = 0.9999999965829791
$ws.Range("E9").Value = 3.417020888463185 / 1000000000

# Row 10
$ws.Range("D10").Value = 0.5395406062475749
$ws.Range("E10").Value = 0.4604593937524251

# Row 11
$ws.Range("D11").Value = 0.9203826970671617
$ws.Range("E11").Value = 0.07961730293283831
$ws.Range("F11").Value = 0.1721367985010147

# Row 12
$ws.Range("D12").Value = 8.63202995278202 / 100000000
$ws.Range("E12").Value = 8.63202995278202 / 100000000

# Row 13
$ws.Range("D13").Value = 0.9999638435713357
$ws.Range("E13").Value = 0.9999638435713357

# Row 14
$ws.Range("D14").Value = 0.001020298220550728
$ws.Range("E14").Value = 0.001020298220550728

# Row 15
$ws.Range("D15").Value = 8.505779214874644 / 100000
$ws.Range("E15").Value = 8.505779214874644 / 100000

# Row 16
$ws.Range("D16").Value = 0.1411269049977135
$ws.Range("E16").Value = 0.1411269049977135

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# Row 18
$ws.Range("D18").Value = 0.9119228338508258
$ws.Range("E18").Value = 0.08807716614917416

# Row 19
$ws.Range("D19").Value = 0.999999999976799
$ws.Range("E19").Value = 2.320099667940667 / 100000000000

# Row 20
$ws.Range("D20").Value = 0.6666683477065343
$ws.Range("E20").Value = 0.3333316522934657

# Row 21
$ws.Range("D21").Value = 0.9625873382108208
$ws.Range("E21").Value = 0.03741266178917924
$ws.Range("F21").Value = 1.091668844223022
